# Nononcology smoke test data
# Rework the scenario2 block to reference the LIVEHTA "NonOncology" project
# (previously the Economic/Takeda rows), drop the old Economic rows and the
# now-unused scenario3/scenario4 blocks, folding their remaining content up
# into the scenario2 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Row 5: was the Economic "StandardExcelReport" row, now carries the
#    LIVEHTA "CompleteExcelReport" filename (previously sitting in I8).
$ws.Cells.Item(5, 9).Value = "CompleteExcelReport-LIVEHTA Automation-Test_NonOncology_Automation_3-Clinical-2023_"

# -- Row 6: was the Economic "ExcelReport" filename, now carries the LIVEHTA
#    "StandardExcelReport" filename (previously sitting in I9).
$ws.Cells.Item(6, 9).Value = "StandardExcelReport-LIVEHTA Automation-Test_NonOncology_Automation_3-Clinical-2023_"

# -- Row 7: was scenario2 / Takeda / Economic header row; becomes the
#    scenario2 / LIVEHTA / Clinical header row carrying what used to be
#    row 9's sub_pop_section1 triple. The old I7 filename is gone.
$ws.Cells.Item(7, 2).Value = "LIVEHTA Automation - Test_NonOncology_Automation_3"
$ws.Cells.Item(7, 3).Value = "LIVEHTA Automation - Test_NonOncology_Automation_3_radio_button"
$ws.Cells.Item(7, 4).Value = "Clinical"
$ws.Cells.Item(7, 5).Value = "Clinical_radio_button"
$ws.Cells.Item(7, 6).Value = "sub_pop_section1"
$ws.Cells.Item(7, 7).Value = "sub_pop_section1_checkbox"
$ws.Cells.Item(7, 8).Value = "sub_pop_section"
$ws.Cells.Item(7, 9).ClearContents()

# -- Row 8: was a spacer row holding only I8; becomes a scenario2 data row
#    carrying what used to be row 10's intervention_section4 triple. It
#    loses the centered "spacer" style that A8 used to carry.
$ws.Cells.Item(8, 1).ClearFormats()
$ws.Cells.Item(8, 1).Value = "scenario2"
$ws.Cells.Item(8, 6).Value = "intervention_section4"
$ws.Cells.Item(8, 7).Value = "intervention_section4_checkbox"
$ws.Cells.Item(8, 8).Value = "intervention_section"
$ws.Cells.Item(8, 9).ClearContents()

# -- Row 9: was scenario3 / LIVEHTA header row; becomes scenario2, carrying
#    what used to be row 11's study_design_section1 triple.
$ws.Cells.Item(9, 1).Value = "scenario2"
$ws.Cells.Item(9, 2).ClearContents()
$ws.Cells.Item(9, 3).ClearContents()
$ws.Cells.Item(9, 4).ClearContents()
$ws.Cells.Item(9, 5).ClearContents()
$ws.Cells.Item(9, 6).Value = "study_design_section1"
$ws.Cells.Item(9, 7).Value = "study_design_section1_checkbox"
$ws.Cells.Item(9, 8).Value = "study_design_section"
$ws.Cells.Item(9, 9).ClearContents()

# -- Row 10: stays scenario3->scenario2, carrying what used to be row 12's
#    reported_variable_section3 triple.
$ws.Cells.Item(10, 1).Value = "scenario2"
$ws.Cells.Item(10, 6).Value = "reported_variable_section3"
$ws.Cells.Item(10, 7).Value = "reported_variable_section3_checkbox"
$ws.Cells.Item(10, 8).Value = "reported_variable_section"

# -- Drop the now-unused trailing rows: old rows 11-12 (rest of scenario3,
#    already folded into rows 9-10 above) and row 14 (scenario4, blank row
#    13 included).
$ws.Rows("11:14").Delete()

# -- Window / view state.
$ws.Range("E1").Select()
$ws.Range("I5:I6").Select()
$aw = $excel.ActiveWindow
$aw.Left = -108
$aw.Top = -108
$aw.Width = 23256
$aw.Height = 12576
